# Apply the "Error Calculations and Plots" missing-data edit.
#
# Summary of the change:
#  - Column D ("C") values are newly imputed/cleared on several rows in the
#    RM block (rows 2-25 in the original layout).
#  - Two entire data rows ("RM 232" and "SC 92") are removed from the bottom
#    block, shifting the remaining SC rows up and changing some of their
#    previously-missing B/D values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column D changes within the RM block (rows 2-25) ---
# Newly filled in:
$ws.Range("D2").Value  = -13.5
$ws.Range("D12").Value = -14.1
$ws.Range("D20").Value = -14
$ws.Range("D21").Value = -14.3

# Newly cleared (became missing):
$ws.Range("D6").ClearContents()
$ws.Range("D14").ClearContents()
$ws.Range("D22").ClearContents()
$ws.Range("D23").ClearContents()

# --- Remove the "RM 232" (row 26) and "SC 92" (row 28) records entirely ---
# Delete the lower row first so the earlier row index stays valid.
$ws.Rows(28).Delete()
$ws.Rows(26).Delete()

# --- Fix up values on the rows that shifted, matching the new data ---
# (Originally "SC 120", now row 30) B value was missing, now filled in.
$ws.Range("B30").Value = -19.7
# (Originally "SC 132", now row 31) D value was missing, now filled in.
$ws.Range("D31").Value = -13.7
# (Originally "SC 193", now row 32) B value was present, now missing.
$ws.Range("B32").ClearContents()
# (Originally "SC 232", now row 33) D value was missing, now filled in.
$ws.Range("D33").Value = -14.1
